$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.23504723507479
$ws.Range("C2").Value = 12.18441315306231
$ws.Range("E2").Value = 14.20124510163991
$ws.Range("F2").Value = 45.95116358346976
$ws.Range("G2").Value = 3.683102720678727
$ws.Range("J2").Value = 9.150573405242156
$ws.Range("M2").Value = 19.85570333019751
$ws.Range("N2").Value = 19.82492251762836

$ws.Range("B3").Value = 18.69473626863942
$ws.Range("C3").Value = 11.66443919940734
$ws.Range("E3").Value = 14.21102626697968
$ws.Range("F3").Value = 45.68811050696194
$ws.Range("G3").Value = 3.687653317454016
$ws.Range("J3").Value = 9.17430836948998
$ws.Range("M3").Value = 19.74759703417283
$ws.Range("N3").Value = 19.87816193193624

$ws.Range("B4").Value = 18.36167746262475
$ws.Range("C4").Value = 11.33734472051807
$ws.Range("E4").Value = 14.21793458570563
$ws.Range("F4").Value = 45.53840437730451
$ws.Range("G4").Value = 3.690588238899446
$ws.Range("J4").Value = 9.190004455167362
$ws.Range("M4").Value = 19.68601777239489
$ws.Range("N4").Value = 19.9128189132649

$ws.Range("B5").Value = 18.22585765160736
$ws.Range("C5").Value = 11.20230926432672
$ws.Range("E5").Value = 14.22097717017351
$ws.Range("F5").Value = 45.48040716430535
$ws.Range("G5").Value = 3.691819808642884
$ws.Range("J5").Value = 9.196682829304251
$ws.Range("M5").Value = 19.66214987038455
$ws.Range("N5").Value = 19.92743606505276

$ws.Range("B6").Value = 18.20330609880418
$ws.Range("C6").Value = 11.17978867274052
$ws.Range("E6").Value = 14.2214961350739
$ws.Range("F6").Value = 45.47095958086675
$ws.Range("G6").Value = 3.692026462347474
$ws.Range("J6").Value = 9.197808799504081
$ws.Range("M6").Value = 19.65826121341107
$ws.Range("N6").Value = 19.92989305870628

$ws.Range("B7").Value = 18.35984582133143
$ws.Range("C7").Value = 11.33553031862875
$ws.Range("E7").Value = 14.21797469781527
$ws.Range("F7").Value = 45.53760997358226
$ws.Range("G7").Value = 3.690604704075166
$ws.Range("J7").Value = 9.190093380156632
$ws.Range("M7").Value = 19.68569089245538
$ws.Range("N7").Value = 19.91301404531733

$ws.Range("B8").Value = 19.04915526590316
$ws.Range("C8").Value = 12.0068788386833
$ws.Range("E8").Value = 14.20443054524045
$ws.Range("F8").Value = 45.85803712665225
$ws.Range("G8").Value = 3.6846426269056
$ws.Range("J8").Value = 9.158524061837587
$ws.Range("M8").Value = 19.81744465228566
$ws.Range("N8").Value = 19.84287065575936

$ws.Range("B9").Value = 20.38088984720274
$ws.Range("C9").Value = 13.25262297492615
$ws.Range("E9").Value = 14.18501600593088
$ws.Range("F9").Value = 46.57812381817424
$ws.Range("G9").Value = 3.674061593408139
$ws.Range("J9").Value = 9.105535547424715
$ws.Range("M9").Value = 20.11294250977691
$ws.Range("N9").Value = 19.72096258689939

$ws.Range("B10").Value = 21.33515210149229
$ws.Range("C10").Value = 14.11484584362959
$ws.Range("E10").Value = 14.17508824862969
$ws.Range("F10").Value = 47.16013813226837
$ws.Range("G10").Value = 3.666955084328369
$ws.Range("J10").Value = 9.07205715580271
$ws.Range("M10").Value = 20.35129978845493
$ws.Range("N10").Value = 19.64097383095274

$ws.Range("B11").Value = 21.76178257736786
$ws.Range("C11").Value = 14.49393307920819
$ws.Range("E11").Value = 14.17150955934355
$ws.Range("F11").Value = 47.43569552904999
$ws.Range("G11").Value = 3.663864990247677
$ws.Range("J11").Value = 9.058014447904366
$ws.Range("M11").Value = 20.4640242433218
$ws.Range("N11").Value = 19.60667403213621

$ws.Range("B12").Value = 21.92209638632984
$ws.Range("C12").Value = 14.63548225336981
$ws.Range("E12").Value = 14.17028889980062
$ws.Range("F12").Value = 47.54152836290922
$ws.Range("G12").Value = 3.662715210364838
$ws.Range("J12").Value = 9.052867839532
$ws.Range("M12").Value = 20.50729818470827
$ws.Range("N12").Value = 19.59398672097009

$ws.Range("B13").Value = 21.88762768551339
$ws.Range("C13").Value = 14.60508780925533
$ws.Range("E13").Value = 14.17054581280391
$ws.Range("F13").Value = 47.51867041602439
$ws.Range("G13").Value = 3.662961932343073
$ws.Range("J13").Value = 9.053968639184559
$ws.Range("M13").Value = 20.49795271949588
$ws.Range("N13").Value = 19.59670574292209

$ws.Range("B14").Value = 21.77499748832469
$ws.Range("C14").Value = 14.50561914456214
$ws.Range("E14").Value = 14.1714064407642
$ws.Range("F14").Value = 47.44437308999306
$ws.Range("G14").Value = 3.66376998972977
$ws.Range("J14").Value = 9.057587602671155
$ws.Range("M14").Value = 20.46757281456607
$ws.Range("N14").Value = 19.60562419487667

$ws.Range("B15").Value = 21.70584178407615
$ws.Range("C15").Value = 14.44442780691272
$ws.Range("E15").Value = 14.17195110916251
$ws.Range("F15").Value = 47.39905520407763
$ws.Range("G15").Value = 3.664267596834284
$ws.Range("J15").Value = 9.059826611625624
$ws.Range("M15").Value = 20.44903986699809
$ws.Range("N15").Value = 19.61112627300858

$ws.Range("B16").Value = 21.30710580972113
$ws.Range("C16").Value = 14.08979708549258
$ws.Range("E16").Value = 14.17534096334843
$ws.Range("F16").Value = 47.14234137884551
$ws.Range("G16").Value = 3.667159887985442
$ws.Range("J16").Value = 9.072998794913108
$ws.Range("M16").Value = 20.34401686560282
$ws.Range("N16").Value = 19.64325750195984

$ws.Range("B17").Value = 21.06045705114099
$ws.Range("C17").Value = 13.86879058963862
$ws.Range("E17").Value = 14.17766042715094
$ws.Range("F17").Value = 46.98757409817156
$ws.Range("G17").Value = 3.668970657039026
$ws.Range("J17").Value = 9.081383747192584
$ws.Range("M17").Value = 20.28066766560258
$ws.Range("N17").Value = 19.66350439282797

$ws.Range("B18").Value = 20.91789733283501
$ws.Range("C18").Value = 13.74044208384138
$ws.Range("E18").Value = 14.17908277992085
$ws.Range("F18").Value = 46.89957648050947
$ws.Range("G18").Value = 3.670025603207444
$ws.Range("J18").Value = 9.086318230332916
$ws.Range("M18").Value = 20.24463696988281
$ws.Range("N18").Value = 19.67534628440628

$ws.Range("B19").Value = 20.86951525630544
$ws.Range("C19").Value = 13.69677771751402
$ws.Range("E19").Value = 14.17957953200748
$ws.Range("F19").Value = 46.8699592550785
$ws.Range("G19").Value = 3.670385102623099
$ws.Range("J19").Value = 9.088008130709909
$ws.Range("M19").Value = 20.23250823449138
$ws.Range("N19").Value = 19.67938944718586

$ws.Range("B20").Value = 21.08678631131873
$ws.Range("C20").Value = 13.8924454286089
$ws.Range("E20").Value = 14.17740438405959
$ws.Range("F20").Value = 47.00394420603101
$ws.Range("G20").Value = 3.66877650769193
$ws.Range("J20").Value = 9.080479593351754
$ws.Range("M20").Value = 20.28736948064458
$ws.Range("N20").Value = 19.66132874067529

$ws.Range("B21").Value = 21.80811467524649
$ws.Range("C21").Value = 14.53489067928106
$ws.Range("E21").Value = 14.17115000514585
$ws.Range("F21").Value = 47.46615625821958
$ws.Range("G21").Value = 3.663532091979484
$ws.Range("J21").Value = 9.056519979551529
$ws.Range("M21").Value = 20.47648042448393
$ws.Range("N21").Value = 19.60299644372536

$ws.Range("B22").Value = 22.27223851202079
$ws.Range("C22").Value = 14.94305021162487
$ws.Range("E22").Value = 14.16784632165847
$ws.Range("F22").Value = 47.77686232251883
$ws.Range("G22").Value = 3.660223236580362
$ws.Range("J22").Value = 9.041858255577639
$ws.Range("M22").Value = 20.6034859593832
$ws.Range("N22").Value = 19.56662954157828

$ws.Range("B23").Value = 22.02524695564501
$ws.Range("C23").Value = 14.72631285085443
$ws.Range("E23").Value = 14.16953792583752
$ws.Range("F23").Value = 47.61026656594899
$ws.Range("G23").Value = 3.661978425004343
$ws.Range("J23").Value = 9.049592106230588
$ws.Range("M23").Value = 20.53539879610369
$ws.Range("N23").Value = 19.58587813380882

$ws.Range("B24").Value = 21.07488520107378
$ws.Range("C24").Value = 13.88175508187344
$ws.Range("E24").Value = 14.17751986433226
$ws.Range("F24").Value = 46.99654022224569
$ws.Range("G24").Value = 3.66886423927125
$ws.Range("J24").Value = 9.080888006677462
$ws.Range("M24").Value = 20.28433837428173
$ws.Range("N24").Value = 19.66231172489482

$ws.Range("B25").Value = 20.02405264372042
$ws.Range("C25").Value = 12.92430148453717
$ws.Range("E25").Value = 14.18950529638118
$ws.Range("F25").Value = 46.37378850960855
$ws.Range("G25").Value = 3.676806147591928
$ws.Range("J25").Value = 9.118914017235017
$ws.Range("M25").Value = 20.02916367507026
$ws.Range("N25").Value = 19.75226322434237

